$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.393.45'
$ws.Range("E2").Value = '  -3.64%  '
$ws.Range("D3").Value = '3.564.08'
$ws.Range("E3").Value = '  -4.15%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '581.71'
$ws.Range("E5").Value = '  -5.12%  '
$ws.Range("D6").Value = '184.64'
$ws.Range("E6").Value = '  -2.30%  '
$ws.Range("D7").Value = '3.560.00'
$ws.Range("E7").Value = '  -4.12%  '
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  -3.83%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").Value = '0.670'
$ws.Range("E10").Value = '  -6.72%  '
$ws.Range("E11").Value = '  -9.83%  '
$ws.Range("D12").Value = '52.90'
$ws.Range("E12").Value = '  -7.46%  '
$ws.Range("D13").Value = '0.0000260'
$ws.Range("E13").Value = '  -10.43%  '
$ws.Range("D14").Value = '9.80'
$ws.Range("E14").Value = '  -7.53%  '
$ws.Range("D15").Value = '4.128.42'
$ws.Range("E15").Value = '  -4.32%  '
$ws.Range("D16").Value = '3.563.91'
$ws.Range("E16").Value = '  -4.25%  '
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '18.35'
$ws.Range("E18").Value = '  -5.36%  '
$ws.Range("D19").Value = '12.19'
$ws.Range("E19").Value = '  -6.44%  '
$ws.Range("D20").Value = '66.203.63'
$ws.Range("E20").Value = '  -3.69%  '
$ws.Range("D21").Value = '1.06'
$ws.Range("E21").Value = '  -7.10%  '
$ws.Range("D22").Value = '395.24'
$ws.Range("E22").Value = '  -3.93%  '
$ws.Range("D23").Value = '4.32'
$ws.Range("E23").Value = '  -5.89%  '
$ws.Range("D24").Value = '85.96'
$ws.Range("E24").Value = '  -3.77%  '
$ws.Range("D25").Value = '11.24'
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").Value = '2.90'
$ws.Range("E26").Value = '  -4.76%  '
$ws.Range("D27").Value = '12.45'
$ws.Range("E27").Value = '  -3.43%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  -6.25%  '
$ws.Range("D30").Value = '8.94'
$ws.Range("E30").Value = '  -7.39%  '
$ws.Range("D31").Value = '31.04'
$ws.Range("E31").Value = '  -6.57%  '
$ws.Range("D32").Value = '7.06'
$ws.Range("E32").Value = '  -3.39%  '
$ws.Range("D33").Value = '12.16'
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("D34").Value = '618.28'
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").Value = '63.63'
$ws.Range("E35").Value = '  -3.87%  '
$ws.Range("E36").Value = '  -8.45%  '
$ws.Range("D37").Value = '41.36'
$ws.Range("E37").Value = '  -7.40%  '
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = '0.396'
$ws.Range("E39").Value = '  -4.76%  '
$ws.Range("D40").Value = '0.0₃0762'
$ws.Range("E40").Value = '  -8.77%  '
$ws.Range("E41").Value = '  -6.21%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '2.986.83'
$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  -7.84%  '
$ws.Range("E45").Value = '  -4.48%  '
$ws.Range("E46").Value = '  -8.08%  '
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").Value = '0.131'
$ws.Range("E48").Value = '  -6.79%  '
$ws.Range("E49").Value = '  -6.85%  '
$ws.Range("D50").Value = '136.96'
$ws.Range("E50").Value = '  -3.67%  '
$ws.Range("D51").Value = '2.72'
$ws.Range("E51").Value = '  -1.98%  '
